$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Ibrahim Al Shamsi"
$summary.Range("B4").Value = 2083.2
$summary.Range("B6").Value = 378298
$summary.Range("B7").Value = 53794
$summary.Range("B8").Value = 324504
$summary.Range("B9").Value = 7.03

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

# Remove the old "Liquid Assets / Savings Account / 11585" row (row 4),
# shifting the TOTAL ASSETS row up from row 5 to row 4.
$assets.Rows.Item(4).Delete()

# Update remaining rows to their new values.
$assets.Range("B2").Value = "Luxury Car"
$assets.Range("C2").Value = 374149

$assets.Range("A3").Value = "Liquid Assets"
$assets.Range("B3").Value = "Savings Account"
$assets.Range("C3").Value = 4149

$assets.Range("C4").Value = 378298

# ---------------------------------------------------------------------------
# Liabilities sheet
# ---------------------------------------------------------------------------
$liabilities = $wb.Worksheets.Item("Liabilities")

# Remove the old "Personal Loans / Personal Loan / ..." row (row 4) and the
# old "Credit Cards / Credit Card Balance / ..." row (row 5), shifting the
# TOTAL LIABILITIES row up from row 6 to row 4.
$liabilities.Rows.Item(5).Delete()
$liabilities.Rows.Item(4).Delete()

# Update remaining rows to their new values.
$liabilities.Range("A2").Value = "Personal Loans"
$liabilities.Range("B2").Value = "Personal Loan"
$liabilities.Range("C2").Value = 27517
$liabilities.Range("D2").Value = 764
$liabilities.Range("E2").Value = 3

$liabilities.Range("A3").Value = "Credit Cards"
$liabilities.Range("B3").Value = "Credit Card Balance"
$liabilities.Range("C3").Value = 26277
$liabilities.Range("D3").Value = 1314
$liabilities.Range("E3").Value = 1

$liabilities.Range("C4").Value = 53794
